$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name filled in (G6) -- copy formatting from the Employee Name
# cell (G4) which already has the "underlined text" look, then set the value.
$ws.Range("G6").Value = "Ankita Gangotra"
$ws.Range("G4").Copy()
$ws.Range("G6").PasteSpecial(-4122)

# Daily hours entered for Wed/Thu/Fri of the second week (rows 15-17, col B)
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 2

# Supervisor sign-off block (row 27): initials + sign-off date
$ws.Range("D27").Value = 41800
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("A27").Value = "A.G"
$ws.Range("A25").Copy()
$ws.Range("A27").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Move the active selection the way it was left when the sheet was saved
$ws.Range("H26").Select()
